# Refresh "想去人数" (interest count, column F) figures to the values
# captured at the newer data-collection run (gh-pages output @ 456a3b4).
# Only column F numeric values change; everything else is untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1016
$ws.Range("F4").Value = 1234
$ws.Range("F5").Value = 1094
$ws.Range("F6").Value = 3251
$ws.Range("F9").Value = 1152
$ws.Range("F10").Value = 711
$ws.Range("F13").Value = 41
$ws.Range("F14").Value = 101
$ws.Range("F16").Value = 1501
$ws.Range("F17").Value = 1501
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 277
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 592
$ws.Range("F22").Value = 317
$ws.Range("F23").Value = 543
$ws.Range("F24").Value = 529
$ws.Range("F25").Value = 42393
$ws.Range("F26").Value = 42394
$ws.Range("F27").Value = 707
$ws.Range("F29").Value = 32120
$ws.Range("F30").Value = 32121
$ws.Range("F31").Value = 422
$ws.Range("F32").Value = 7
$ws.Range("F35").Value = 914
$ws.Range("F36").Value = 218
$ws.Range("F37").Value = 145
$ws.Range("F38").Value = 478
$ws.Range("F39").Value = 1138
$ws.Range("F40").Value = 5285
$ws.Range("F41").Value = 678
$ws.Range("F42").Value = 404
$ws.Range("F43").Value = 12
$ws.Range("F45").Value = 307
$ws.Range("F48").Value = 13
$ws.Range("F49").Value = 36

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 212
$ws.Range("F13").Value = 1788
$ws.Range("F15").Value = 715
$ws.Range("F18").Value = 388
$ws.Range("F20").Value = 57
$ws.Range("F35").Value = 1116
$ws.Range("F38").Value = 75
$ws.Range("F39").Value = 75
$ws.Range("F45").Value = 793

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 686
$ws.Range("F5").Value = 517
$ws.Range("F6").Value = 512

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 686
$ws.Range("F4").Value = 517
$ws.Range("F6").Value = 1016
$ws.Range("F7").Value = 1234
$ws.Range("F8").Value = 212
$ws.Range("F9").Value = 1094
$ws.Range("F10").Value = 3251
$ws.Range("F13").Value = 1152
$ws.Range("F14").Value = 711
$ws.Range("F15").Value = 512
$ws.Range("F17").Value = 1788
$ws.Range("F20").Value = 41
$ws.Range("F21").Value = 101
$ws.Range("F23").Value = 1501
$ws.Range("F24").Value = 1501
$ws.Range("F25").Value = 277
$ws.Range("F27").Value = 22
$ws.Range("F28").Value = 592
$ws.Range("F29").Value = 388
$ws.Range("F30").Value = 317
$ws.Range("F31").Value = 529
$ws.Range("F32").Value = 42401
$ws.Range("F34").Value = 707
$ws.Range("F36").Value = 32122
$ws.Range("F37").Value = 422
$ws.Range("F38").Value = 914
$ws.Range("F39").Value = 218
$ws.Range("F40").Value = 145
$ws.Range("F41").Value = 478
$ws.Range("F42").Value = 1138
$ws.Range("F43").Value = 5285
$ws.Range("F44").Value = 678
$ws.Range("F46").Value = 404
$ws.Range("F47").Value = 75
$ws.Range("F48").Value = 12
$ws.Range("F49").Value = 307
$ws.Range("F51").Value = 793
$ws.Range("F52").Value = 36
